$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 160-161 (pushing the existing weekly records down),
# adding a new "Betarraga" price week (2021-11-11 / serial 44511) at the
# top of this block.
$ws.Rows("160:161").Insert()

# Row 160 - Primera
$ws.Range("A160").Value = 8
$ws.Range("B160").Value = "Terminal La Palmera de La Serena"
$ws.Range("C160").Value = "Coquimbo"
$ws.Range("D160").Value = 44511
$ws.Range("E160").Value = 4
$ws.Range("F160").Value = 100114014
$ws.Range("G160").Value = "Betarraga"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 3000
$ws.Range("K160").Value = 450
$ws.Range("L160").Value = 500
$ws.Range("M160").Value = 475
$ws.Range("N160").Value = "`$/paquete 3 unidades"
$ws.Range("O160").Value = "Provincia del Elquí"
$ws.Range("P160").Value = 158
$ws.Range("Q160").Value = 3
$ws.Range("R160").Value = "Hortaliza"

# Row 161 - Segunda
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 44511
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100114014
$ws.Range("G161").Value = "Betarraga"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 1500
$ws.Range("K161").Value = 350
$ws.Range("L161").Value = 400
$ws.Range("M161").Value = 375
$ws.Range("N161").Value = "`$/paquete 3 unidades"
$ws.Range("O161").Value = "Provincia del Elquí"
$ws.Range("P161").Value = 125
$ws.Range("Q161").Value = 3
$ws.Range("R161").Value = "Hortaliza"
